$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set floodmedia column (H) for data rows 2-27 to "None"
$ws.Range("H2:H27").Value = "None"

# Update selection to match the new active range (H2:H27)
$ws.Range("H2:H27").Select()
